# Error Calculations and Plots
# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# which becomes row 27 after the first delete), then fill in / clear a
# couple of previously-missing "C" column (column D) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - everything below shifts up one row.
$ws.Rows.Item(26).Delete()

# After that delete, the row that was "SC 92" (originally row 28) is now
# row 27. Delete it too - everything below shifts up one more row.
$ws.Rows.Item(27).Delete()

# Now the data has shifted up by two rows total, giving us (in new
# numbering): 26=SC 5, 27=SC 101, 28=SC 105, 29=SC 119, 30=SC 120,
# 31=SC 132, 32=SC 193, 33=SC 232.
# Update the newly-imputed "C" column (column D) values for the rows that
# had missing data.
$ws.Range("D26").Value = -13.8    # SC 5
$ws.Range("D27").Value = ""       # SC 101 (now missing)
$ws.Range("D30").Value = -13.6    # SC 120
$ws.Range("D32").Value = ""       # SC 193 (now missing)
